$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value. These cells hold text-formatted
# numbers/percentages (t="inlineStr" in the source sheet), so each
# write temporarily forces Text number-format to stop Excel from
# auto-converting the numeric-looking string into a real number,
# then restores the original (default) style afterwards.
$updates = @{
    "D2" = "331.35"
    "E2" = "0.76%"
    "D3" = "41.24"
    "E3" = "2.16%"
    "D4" = "5.686"
    "E4" = "-4.49%"
    "D5" = "0.08073"
    "E5" = "-0.71%"
    "D6" = "2.035"
    "E6" = "3.44%"
    "D7" = "8.753"
    "E7" = "-0.09%"
    "D8" = "4.536"
    "E8" = "-1.56%"
    "E9" = "0.92%"
    "D10" = "0.9241"
    "E10" = "-2.49%"
    "D11" = "0.1262"
    "E11" = "-4.17%"
    "D12" = "0.1949"
    "E12" = "-2.08%"
    "D13" = "8.734"
    "E13" = "-3.45%"
    "D14" = "0.09388"
    "E14" = "0.92%"
    "D15" = "0.03748"
    "E15" = "7.85%"
    "D16" = "0.1053"
    "E16" = "9.53%"
    "D17" = "0.001300"
    "E17" = "-1.24%"
    "D18" = "0.006308"
    "E18" = "-1.42%"
    "D19" = "3.364"
    "E19" = "0.28%"
    "E20" = "-1.73%"
    "D21" = "0.1418"
    "E21" = "1.08%"
    "D22" = "0.2654"
    "D23" = "0.04440"
    "E23" = "-0.05%"
    "D24" = "0.001261"
    "E24" = "-0.13%"
    "E25" = "-3.66%"
    "D26" = "0.0001243"
    "E26" = "13.41%"
    "E39" = "15.49%"
    "D40" = "0.05484"
    "E40" = "3.91%"
    "D41" = "0.007783"
    "E41" = "2.97%"
    "D42" = "0.009942"
    "E42" = "10.19%"
    "D43" = "0.1422"
    "E43" = "-0.73%"
    "D44" = "0.002125"
    "E44" = "3.09%"
    "D45" = "0.01187"
    "E45" = "12.89%"
    "D46" = "0.00006756"
    "E46" = "-0.93%"
    "E47" = "-0.05%"
    "D48" = "0.002282"
    "E48" = "26.55%"
    "D49" = "0.003014"
    "E49" = "-14.22%"
    "D50" = "0.00002103"
    "E50" = "-0.05%"
    "D51" = "0.0002003"
    "E51" = "-0.05%"
}

foreach ($cell in $updates.Keys) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$cell]
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}
